$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 'Datos actualizados a 20 de Marzo de 2020 a las 03:16'
$ws.Cells.Item(9, 2).Value = 14299
$ws.Cells.Item(9, 3).Value = 5040
$ws.Cells.Item(9, 5).Value = 13960
$ws.Cells.Item(41, 2).Value = 309
$ws.Cells.Item(41, 3).Value = 82
$ws.Cells.Item(41, 5).Value = 269
$ws.Cells.Item(116, 1).Value = 'Camerun'
$ws.Cells.Item(116, 3).Value = 0
$ws.Cells.Item(117, 1).Value = 'Paraguay'
$ws.Cells.Item(117, 3).Value = 2
$ws.Cells.Item(117, 6).Value = 1
$ws.Cells.Item(118, 1).Value = 'Maldivas'
$ws.Cells.Item(118, 3).Value = 0
$ws.Cells.Item(118, 6).Value = 0
$ws.Cells.Item(119, 1).Value = 'Montenegro'
$ws.Cells.Item(119, 3).Value = 5
$ws.Cells.Item(130, 1).Value = 'Etiopia'
$ws.Cells.Item(130, 3).Value = 1
$ws.Cells.Item(131, 1).Value = 'Kenia'
$ws.Cells.Item(131, 3).Value = 0
$ws.Cells.Item(133, 1).Value = 'Puerto Rico'
$ws.Cells.Item(133, 3).Value = 0
$ws.Cells.Item(137, 1).Value = 'Polinesia Francesa'
$ws.Cells.Item(137, 3).Value = 1
$ws.Cells.Item(140, 1).Value = 'Guyana'
$ws.Cells.Item(140, 4).Value = 0
$ws.Cells.Item(140, 8).Value = 1
$ws.Cells.Item(141, 1).Value = 'Aruba'
$ws.Cells.Item(141, 4).Value = 1
$ws.Cells.Item(141, 8).Value = 0
$ws.Cells.Item(143, 1).Value = 'San Bartolome'
$ws.Cells.Item(143, 3).Value = 0
$ws.Cells.Item(144, 1).Value = 'San Martin (Parte Francesa)'
$ws.Cells.Item(144, 3).Value = 0
$ws.Cells.Item(145, 1).Value = 'Kirguistan'
$ws.Cells.Item(145, 3).Value = 0
$ws.Cells.Item(146, 1).Value = 'Bahamas'
$ws.Cells.Item(146, 3).Value = 2
$ws.Cells.Item(147, 1).Value = 'Islas Virgenes de los Estados Unidos'
$ws.Cells.Item(147, 3).Value = 1
$ws.Cells.Item(148, 1).Value = 'Gabon'
$ws.Cells.Item(149, 1).Value = 'Namibia'
$ws.Cells.Item(149, 3).Value = 1
$ws.Cells.Item(155, 1).Value = 'Santa Lucia'
$ws.Cells.Item(156, 1).Value = 'Zambia'
$ws.Cells.Item(159, 1).Value = 'Mauritania'
$ws.Cells.Item(160, 1).Value = 'Groenlandia'
$ws.Cells.Item(162, 1).Value = 'Isla de Man'
$ws.Cells.Item(163, 1).Value = 'Guinea'
$ws.Cells.Item(164, 1).Value = 'San Martin (Parte Holandesa)'
$ws.Cells.Item(165, 1).Value = 'Togo'
$ws.Cells.Item(166, 1).Value = 'Santa Sede'
$ws.Cells.Item(167, 1).Value = 'Surinam'
$ws.Cells.Item(168, 1).Value = 'Somalia'
$ws.Cells.Item(168, 3).Value = 0
$ws.Cells.Item(169, 1).Value = 'San Vicente y las Granadinas'
$ws.Cells.Item(170, 1).Value = 'Fiyi'
$ws.Cells.Item(171, 1).Value = 'El Salvador'
$ws.Cells.Item(171, 3).Value = 1
$ws.Cells.Item(172, 1).Value = 'Republica de Yibuti'
$ws.Cells.Item(173, 1).Value = 'Montserrat'
$ws.Cells.Item(174, 1).Value = 'Gambia'
$ws.Cells.Item(175, 1).Value = 'Nicaragua'
$ws.Cells.Item(177, 1).Value = 'Suazilandia'
$ws.Cells.Item(179, 1).Value = 'Republica del Chad'
$ws.Cells.Item(180, 1).Value = 'Butan'
$ws.Cells.Item(181, 1).Value = 'Republica de Africa Central'
